$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows after the current last row (212), shifting nothing
# else, to make room for rows 213-216.
$ws.Range("A213:T216").Insert()

# The 4 new rows preserve the original (pre-edit) data that used to live in
# rows 209-212 (Fukumoto entries dated 44399/44399/44400/44400), so copy it
# down before rows 209-212 get overwritten with the new values below.
$ws.Range("A209:T212").Copy($ws.Range("A213:T216"))

# --- Row 209: Fukumoto -> Cara cara, new date/volume/price data ---
$ws.Cells.Item(209, 4).Value = 44448
$ws.Cells.Item(209, 11).Value = "Cara cara"
$ws.Cells.Item(209, 13).Value = 240
$ws.Cells.Item(209, 14).Value = 5500
$ws.Cells.Item(209, 15).Value = 6000
$ws.Cells.Item(209, 16).Value = 5750
$ws.Cells.Item(209, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(209, 19).Value = 359
$ws.Cells.Item(209, 20).Value = 16

# --- Row 210: Fukumoto -> Cara cara, new date/volume/price data ---
$ws.Cells.Item(210, 4).Value = 44448
$ws.Cells.Item(210, 11).Value = "Cara cara"
$ws.Cells.Item(210, 13).Value = 180
$ws.Cells.Item(210, 14).Value = 4500
$ws.Cells.Item(210, 15).Value = 5000
$ws.Cells.Item(210, 16).Value = 4750
$ws.Cells.Item(210, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(210, 19).Value = 297
$ws.Cells.Item(210, 20).Value = 16

# --- Row 211: Fukumoto -> Navel Late, new date/volume/price data ---
$ws.Cells.Item(211, 4).Value = 44448
$ws.Cells.Item(211, 11).Value = "Navel Late"
$ws.Cells.Item(211, 13).Value = 240
$ws.Cells.Item(211, 14).Value = 5500
$ws.Cells.Item(211, 15).Value = 6000
$ws.Cells.Item(211, 16).Value = 5750
$ws.Cells.Item(211, 19).Value = 383

# --- Row 212: now Navel Late / Segunda with new date/volume/price data ---
$ws.Cells.Item(212, 4).Value = 44448
$ws.Cells.Item(212, 5).Value = 16
$ws.Cells.Item(212, 6).Value = "Fruta"
$ws.Cells.Item(212, 7).Value = 100102
$ws.Cells.Item(212, 8).Value = "Cítricos"
$ws.Cells.Item(212, 9).Value = 100102005
$ws.Cells.Item(212, 10).Value = "Naranja"
$ws.Cells.Item(212, 11).Value = "Navel Late"
$ws.Cells.Item(212, 12).Value = "Segunda"
$ws.Cells.Item(212, 13).Value = 180
$ws.Cells.Item(212, 14).Value = 4500
$ws.Cells.Item(212, 15).Value = 5000
$ws.Cells.Item(212, 16).Value = 4750
$ws.Cells.Item(212, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(212, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(212, 19).Value = 317
$ws.Cells.Item(212, 20).Value = 15
